$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R4").Copy() | Out-Null
$ws.Range("S4").PasteSpecial(-4122) | Out-Null
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy() | Out-Null
$ws.Range("S5").PasteSpecial(-4122) | Out-Null
$ws.Range("S5").Value = 4.9000000000000004

$ws.Range("R6").Copy() | Out-Null
$ws.Range("S6").PasteSpecial(-4122) | Out-Null
$ws.Range("S6").Value = 6.1

$ws.Range("R7").Copy() | Out-Null
$ws.Range("S7").PasteSpecial(-4122) | Out-Null
$ws.Range("S7").Value = 4

$ws.Range("R8").Copy() | Out-Null
$ws.Range("S8").PasteSpecial(-4122) | Out-Null
$ws.Range("S8").Font.Bold = $true
$ws.Range("S8").Font.Italic = $true

$ws.Range("R9").Copy() | Out-Null
$ws.Range("S9").PasteSpecial(-4122) | Out-Null
$ws.Range("S9").Value = 6.1

$ws.Range("R10").Copy() | Out-Null
$ws.Range("S10").PasteSpecial(-4122) | Out-Null
$ws.Range("S10").Value = 12.4

$ws.Range("R11").Copy() | Out-Null
$ws.Range("S11").PasteSpecial(-4122) | Out-Null
$ws.Range("S11").Value = 3.2

$ws.Range("R12").Copy() | Out-Null
$ws.Range("S12").PasteSpecial(-4122) | Out-Null
$ws.Range("S12").Value = 10.8

$ws.Range("R13").Copy() | Out-Null
$ws.Range("S13").PasteSpecial(-4122) | Out-Null
$ws.Range("S13").Value = 14.6

$ws.Range("R14").Copy() | Out-Null
$ws.Range("S14").PasteSpecial(-4122) | Out-Null
$ws.Range("S14").Value = 8.5

$ws.Range("R15").Copy() | Out-Null
$ws.Range("S15").PasteSpecial(-4122) | Out-Null
$ws.Range("S15").Value = 5.5

$ws.Range("R16").Copy() | Out-Null
$ws.Range("S16").PasteSpecial(-4122) | Out-Null
$ws.Range("S16").Value = 7.1

$ws.Range("R17").Copy() | Out-Null
$ws.Range("S17").PasteSpecial(-4122) | Out-Null
$ws.Range("S17").Value = 4.4000000000000004

$ws.Range("R18").Copy() | Out-Null
$ws.Range("S18").PasteSpecial(-4122) | Out-Null
$ws.Range("S18").Value = 5.8

$ws.Range("R19").Copy() | Out-Null
$ws.Range("S19").PasteSpecial(-4122) | Out-Null
$ws.Range("S19").Value = 11.6

$ws.Range("R20").Copy() | Out-Null
$ws.Range("S20").PasteSpecial(-4122) | Out-Null
$ws.Range("S20").Value = 3.1

$ws.Range("R21").Copy() | Out-Null
$ws.Range("S21").PasteSpecial(-4122) | Out-Null
$ws.Range("S21").Value = 1.5

$ws.Range("R22").Copy() | Out-Null
$ws.Range("S22").PasteSpecial(-4122) | Out-Null
$ws.Range("S22").Value = 2.2999999999999998

$ws.Range("R23").Copy() | Out-Null
$ws.Range("S23").PasteSpecial(-4122) | Out-Null
$ws.Range("S23").Value = 1

$ws.Range("R24").Copy() | Out-Null
$ws.Range("S24").PasteSpecial(-4122) | Out-Null
$ws.Range("S24").Value = 2.2999999999999998

$ws.Range("R25").Copy() | Out-Null
$ws.Range("S25").PasteSpecial(-4122) | Out-Null
$ws.Range("S25").Value = 3.3

$ws.Range("R26").Copy() | Out-Null
$ws.Range("S26").PasteSpecial(-4122) | Out-Null
$ws.Range("S26").Value = 1.6

$ws.Range("R27").Copy() | Out-Null
$ws.Range("S27").PasteSpecial(-4122) | Out-Null
$ws.Range("S27").Value = 4.5999999999999996

$ws.Range("R28").Copy() | Out-Null
$ws.Range("S28").PasteSpecial(-4122) | Out-Null
$ws.Range("S28").Value = 4.4000000000000004

$ws.Range("R29").Copy() | Out-Null
$ws.Range("S29").PasteSpecial(-4122) | Out-Null
$ws.Range("S29").Value = 4.7

$ws.Range("R30").Copy() | Out-Null
$ws.Range("S30").PasteSpecial(-4122) | Out-Null
$ws.Range("S30").Value = 4

$ws.Range("R31").Copy() | Out-Null
$ws.Range("S31").PasteSpecial(-4122) | Out-Null
$ws.Range("S31").Value = 3.2

$ws.Range("R32").Copy() | Out-Null
$ws.Range("S32").PasteSpecial(-4122) | Out-Null
$ws.Range("S32").Value = 4.7

$ws.Range("R33").Copy() | Out-Null
$ws.Range("S33").PasteSpecial(-4122) | Out-Null
$ws.Range("S33").Value = 2.6

$ws.Range("R34").Copy() | Out-Null
$ws.Range("S34").PasteSpecial(-4122) | Out-Null
$ws.Range("S34").Value = 3.3

$ws.Range("R35").Copy() | Out-Null
$ws.Range("S35").PasteSpecial(-4122) | Out-Null
$ws.Range("S35").Value = 2.2000000000000002

$ws.Range("R36").Copy() | Out-Null
$ws.Range("S36").PasteSpecial(-4122) | Out-Null
$ws.Range("S36").Font.Bold = $true
$ws.Range("S36").Font.Italic = $true

$ws.Range("R37").Copy() | Out-Null
$ws.Range("S37").PasteSpecial(-4122) | Out-Null
$ws.Range("S37").Value = 13.2

$ws.Range("R38").Copy() | Out-Null
$ws.Range("S38").PasteSpecial(-4122) | Out-Null
$ws.Range("S38").Value = 7.5

$ws.Range("R39").Copy() | Out-Null
$ws.Range("S39").PasteSpecial(-4122) | Out-Null
$ws.Range("S39").Value = 4.0999999999999996

$ws.Range("R40").Copy() | Out-Null
$ws.Range("S40").PasteSpecial(-4122) | Out-Null
$ws.Range("S40").Value = 4.3

$ws.Range("R41").Copy() | Out-Null
$ws.Range("S41").PasteSpecial(-4122) | Out-Null
$ws.Range("S41").Value = 2.6

$ws.Range("R42").Copy() | Out-Null
$ws.Range("S42").PasteSpecial(-4122) | Out-Null
$ws.Range("S42").Value = 1

$ws.Range("R43").Copy() | Out-Null
$ws.Range("S43").PasteSpecial(-4122) | Out-Null
$ws.Range("S43").Value = "…"

$ws.Range("T12").Select() | Out-Null